$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Remove the old sample rows (3-4) below the header, leaving row 2 for new data
$ws.Range("A3:B4").EntireRow.Delete() | Out-Null

# Header row
$headers = @("name", "msv", "class", "hdcm.uv1", "hdcm.uv2", "hdcm.uv3", "hdcm.uv4", "hdcm.uv5", "hd.01", "hd.02", "hd.03", "pb")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("C1:L1").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral

# Approximate "best fit" column widths for the new columns
$colWidths = @(27.666666666666668, 22.166666666666668, 22.666666666666668, 14.333333333333334, 34.833333333333336, 38.333333333333336, 12.666666666666666, 25.666666666666668, 56.0, 19.833333333333332, 12.666666666666666, 12.666666666666666)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i]
}

# Turn the header row into a table (data rows are appended below afterwards)
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:L1"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleLight1"

# Data row
$ws.Range("A2").Value = "Nguyễn Tiến Binh"
$ws.Range("B2").Value = "B20DCDT021"
$ws.Range("C2").Value = "D20DTMT1"
$ws.Range("L2").Value = "Trần Thị Thúy Hà - C2.3: 1 - C3.2: 2 - C4.1: 3 - C6.1: 5 - C6.2: 6 - GPA: 7"
$ws.Range("A2:B2").Style = "Normal"
